# "moves curves into params.xlsx/params.json"
# Adds a new "curves" worksheet as the first sheet in the workbook and
# populates it with the month / curve-shape lookup table.

$wb = $excel.ActiveWorkbook

# --- create the new sheet and put it first -------------------------------
$curves = $wb.Worksheets.Add()
$curves.Name = "curves"
$curves.Move($wb.Worksheets.Item(1))
$curves.Activate()

# --- header row ------------------------------------------------------------
# Written in this order so new shared-string entries land at the same
# indices the source workbook used (B, D, E, then A).
$curves.Range("B1").Value = "Sudden shock"
$curves.Range("D1").Value = "Shallow mid-term"
$curves.Range("E1").Value = "Sustained impact"
$curves.Range("A1").Value = "month"
$curves.Range("C1").Value = "Follow the curve"

# --- data rows --------------------------------------------------------------
# Each entry: month, Sudden shock, Follow the curve, Shallow mid-term, Sustained impact
$rows = @(
    @(0,  0,      0,      0,      0),
    @(1,  0.0612, 0.0119, 0.0076, 0.0118),
    @(2,  0.1224, 0.0238, 0.0152, 0.0234),
    @(3,  0.2449, 0.0475, 0.0227, 0.0374),
    @(4,  0.2041, 0.095,  0.0303, 0.042),
    @(5,  0.1633, 0.1188, 0.0379, 0.0448),
    @(6,  0.102,  0.1206, 0.0455, 0.0467),
    @(7,  0.051,  0.1069, 0.053,  0.0467),
    @(8,  0.0306, 0.0841, 0.0606, 0.0467),
    @(9,  0.0153, 0.0722, 0.0682, 0.0467),
    @(10, 0.0052, 0.0594, 0.0758, 0.0467),
    @(11, 0,      0.0487, 0.0832, 0.0467),
    @(12, 0,      0.0392, 0.0832, 0.0467),
    @(13, 0,      0.0333, 0.0758, 0.0467),
    @(14, 0,      0.0285, 0.0682, 0.0467),
    @(15, 0,      0.0249, 0.0606, 0.0467),
    @(16, 0,      0.0214, 0.053,  0.0467),
    @(17, 0,      0.0178, 0.0455, 0.0467),
    @(18, 0,      0.0143, 0.0379, 0.0467),
    @(19, 0,      0.0116, 0.0303, 0.0467),
    @(20, 0,      0.009,  0.0227, 0.0467),
    @(21, 0,      0.0062, 0.0152, 0.0467),
    @(22, 0,      0.0034, 0.0076, 0.0467),
    @(23, 0,      0.0015, 0,      0.0467)
)

$r = 2
foreach ($entry in $rows) {
    $curves.Cells.Item($r, 1).Value = [double]$entry[0]
    $curves.Cells.Item($r, 2).Value = [double]$entry[1]
    $curves.Cells.Item($r, 3).Value = [double]$entry[2]
    $curves.Cells.Item($r, 4).Value = [double]$entry[3]
    $curves.Cells.Item($r, 5).Value = [double]$entry[4]
    $r++
}

# --- column widths (approximate Excel's "best fit" autosize for these
# headers/values; the COM ColumnWidth property is in character units and
# gets re-quantized on save, so these are tuned to land as close as
# possible to the original bestFit widths of 6.85, 13.29, 15.86, 17, 16.29)
$curves.Columns("A").ColumnWidth = 6.02
$curves.Columns("B").ColumnWidth = 12.45
$curves.Columns("C").ColumnWidth = 15.02
$curves.Columns("D").ColumnWidth = 16.17
$curves.Columns("E").ColumnWidth = 15.45
